$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add a new row (52) for "Longest Increasing Path in a Matrix" (leetcode 329),
# following the same pattern as the existing backtracking rows.
$ws.Range("A52").Value = 51
$ws.Range("B52").Value = "Longest Increasing Path in a Matrix"
$ws.Range("C52").Value = "backtracking"
$ws.Range("D52").Value = "Array"
$ws.Range("E52").Value = "medium"
$ws.Range("F52").Value = "leetcode 329"

# Match styles used by the neighboring rows (center alignment on A/D/E/F,
# left alignment on B; C keeps the default style).
$ws.Range("A52").HorizontalAlignment = -4108
$ws.Range("B52").HorizontalAlignment = -4131
$ws.Range("D52").HorizontalAlignment = -4108
$ws.Range("E52").HorizontalAlignment = -4108
$ws.Range("F52").HorizontalAlignment = -4108

# Move the active selection to F52, matching the saved view state.
$ws.Range("F52").Select()
